$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# New data filled into column A (네이버 / Naver) for rows 2-4
$ws.Range("A2").Value = "모두비허용"
$ws.Range("A3").Value = "User-agent: *"
$ws.Range("A4").Value = "Disallow: /"

# Existing summary row additions for previously-empty columns
$ws.Range("D2").Value = "존재x"
$ws.Range("J2").Value = "존재x"
$ws.Range("M2").Value = "존재x"

# Correct the Bomtoon (F) summary value
$ws.Range("F2").Value = "모두 허용"

# New data filled into column G (뿌딩 / Pudding) for rows 2-4
$ws.Range("G2").Value = "모두 허용"
$ws.Range("G3").Value = "User-agent: * "
$ws.Range("G4").Value = "Allow : /"

# Comica (I) summary value
$ws.Range("I2").Value = "됨"

# Move/restore the active selection to F3
$ws.Range("F3").Select()
